$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking price strings
# (e.g. "596.79") are not auto-coerced into numbers by the COM layer,
# matching the source inlineStr cell type. Style is restored to Normal
# afterwards so no cell ends up with a different style than before.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '67.719.86'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.782.13'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '596.79'
$ws.Range('D6').Value = '168.98'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').Value = '3.781.58'
$ws.Range('E7').Value = '  -1.82%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -0.83%  '
$ws.Range('D10').Value = '0.164'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('D11').Value = '6.50'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').Value = '0.453'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('E13').Value = '  +3.78%  '
$ws.Range('D14').Value = '36.41'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').Value = '4.416.27'
$ws.Range('E15').Value = '  -1.80%  '
$ws.Range('D16').Value = '3.787.43'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').Value = '18.56'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '67.655.30'
$ws.Range('D19').Value = '7.19'
$ws.Range('E19').Value = '  -2.47%  '
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').Value = '10.52'
$ws.Range('E21').Value = '  -6.21%  '
$ws.Range('D22').Value = '468.40'
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('E23').Value = '  -1.88%  '
$ws.Range('E24').Value = '  -7.74%  '
$ws.Range('D25').Value = '83.80'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = '2.21'
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').Value = '12.14'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '10.33'
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  -1.61%  '
$ws.Range('D31').Value = '3.933.66'
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('E32').Value = '  -0.82%  '
$ws.Range('D33').Value = '30.54'
$ws.Range('E33').Value = '  -2.94%  '
$ws.Range('E34').Value = '  -3.65%  '
$ws.Range('E35').Value = '  -2.31%  '
$ws.Range('D36').Value = '3.747.55'
$ws.Range('E36').Value = '  -1.84%  '
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('E40').Value = '  -2.06%  '
$ws.Range('E41').Value = '  -2.29%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '0.311'
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('D47').Value = '45.73'
$ws.Range('E47').Value = '  -2.63%  '
$ws.Range('D48').Value = '395.86'
$ws.Range('E48').Value = '  -5.39%  '
$ws.Range('E49').Value = '  -9.11%  '
$ws.Range('D50').Value = '140.45'
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('D51').Value = '39.35'
$ws.Range('E51').Value = '  +3.18%  '

$ws.Range('D2:D51').Style = 'Normal'
